$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.457.35'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '2.599.25'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Formula = '="523.16"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").Formula = '="143.81"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("D7").Formula = '="0.998"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Formula = '="0.570"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +0.53%  '
$ws.Range("D9").Value = '2.620.95'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").Formula = '="6.65"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("D11").Formula = '="0.102"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -1.46%  '
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").Value = '3.054.18'
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("D15").Value = '58.250.86'
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("D16").Formula = '="20.54"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  -2.12%  '
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '2.596.15'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").Formula = '="339.53"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("D20").Formula = '="4.37"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").Formula = '="10.30"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("E22").Value = '  +1.82%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Formula = '="65.36"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +1.46%  '
$ws.Range("E25").Value = '  +0.51%  '
$ws.Range("D26").Formula = '="0.404"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -2.81%  '
$ws.Range("D27").Value = '2.717.21'
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("D28").Formula = '="0.996"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("D30").Value = '0.0₃0750'
$ws.Range("E30").Value = '  -5.28%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Formula = '="6.24"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  -6.12%  '
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("D35").Formula = '="149.86"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").Formula = '="4.04"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -1.63%  '
$ws.Range("E37").Value = '  -3.70%  '
$ws.Range("D38").Formula = '="0.874"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -1.74%  '
$ws.Range("D39").Formula = '="0.863"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +1.37%  '
$ws.Range("D40").Formula = '="36.05"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("E42").Value = '  -2.06%  '
$ws.Range("D43").Formula = '="0.997"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").Formula = '="273.61"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("D45").Formula = '="0.599"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("D46").Formula = '="0.0960"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("D47").Formula = '="10.68"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  +0.66%  '
$ws.Range("D48").Formula = '="18.86"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("D50").Formula = '="19.00"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +4.10%  '
$ws.Range("D51").Value = '1.976.50'
$ws.Range("E51").Value = '  -3.00%  '

$excel.CutCopyMode = $false
